$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Remove the stray _GoBack bookmark sitting at the very start of
#    the document (empty paragraph before the cover-page sdt block).
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2. Typo / wording fixes throughout the body text.
# ------------------------------------------------------------------

# "...we Will explain..." -> "...we will explain..."
$r = $d.Content
$r.Find.Execute("we Will explain", $true, $false, $false, $false, $false, $true, 1, $false, "we will explain", 2) | Out-Null

# "...see AnnexI" -> "...see Annex" (only touch the "AnnexI" run so the
# neighbouring "see" run/formatting is left completely untouched)
$r = $d.Content
$r.Find.Execute("AnnexI", $true, $false, $false, $false, $false, $true, 1, $false, "Annex", 2) | Out-Null

# "...in a máximum of..." -> "...in a maximum of..."
$r = $d.Content
$r.Find.Execute(", in a máximum of", $true, $false, $false, $false, $false, $true, 1, $false, ", in a maximum of", 2) | Out-Null

# "...the dayly amortization..." -> "...the daily amortization..."
$r = $d.Content
$r.Find.Execute("the dayly amortization", $true, $false, $false, $false, $false, $true, 1, $false, "the daily amortization", 2) | Out-Null

# "...amortization Will be..." -> "...amortization will be..."
$r = $d.Content
$r.Find.Execute("amortization Will be", $true, $false, $false, $false, $false, $true, 1, $false, "amortization will be", 2) | Out-Null

# "Ristks" -> "Risks"
$r = $d.Content
$r.Find.Execute("Ristks", $true, $false, $false, $false, $false, $true, 1, $false, "Risks", 2) | Out-Null

# "Becaouse" -> "Because"
$r = $d.Content
$r.Find.Execute("Becaouse", $true, $false, $false, $false, $false, $true, 1, $false, "Because", 2) | Out-Null

# "...a smoll Project..." -> "...a small Project..."
$r = $d.Content
$r.Find.Execute("a smoll Project", $true, $false, $false, $false, $false, $true, 1, $false, "a small Project", 2) | Out-Null

# "...are no ristks of..." -> "...are no risks of..."
$r = $d.Content
$r.Find.Execute("are no ristks of", $true, $false, $false, $false, $false, $true, 1, $false, "are no risks of", 2) | Out-Null

# "...earnings Will be 20%..." -> "...earnings will be 20%..."
$r = $d.Content
$r.Find.Execute("earnings Will be 20%", $true, $false, $false, $false, $false, $true, 1, $false, "earnings will be 20%", 2) | Out-Null

# ------------------------------------------------------------------
# 3. "The sum of directa n indirect costs" -> "The sum of direct and
#    indirect costs", and re-insert the _GoBack bookmark right after
#    the newly typed "and" (matching where Word last left the cursor
#    after this edit).
# ------------------------------------------------------------------
$r = $d.Content
$found = $r.Find.Execute(" directa n ", $true, $false, $false, $false, $false, $true, 1, $false, " direct and ", 2)
if ($found) {
    $bmStart = $r.Start + 11
    $bmRange = $d.Range($bmStart, $bmStart)
    $d.Bookmarks.Add("_GoBack", $bmRange)
}
